$d = $word.ActiveDocument

$d.Content.Find.Execute("44+39=83", $true, $false, $false, $false, $false, $true, 1, $false, "39+45=84", 2) | Out-Null
$d.Content.Find.Execute("65+10=75", $true, $false, $false, $false, $false, $true, 1, $false, "9-5=4", 2) | Out-Null
$d.Content.Find.Execute("56-27=29", $true, $false, $false, $false, $false, $true, 1, $false, "13-0=13", 2) | Out-Null
$d.Content.Find.Execute("77+0=77", $true, $false, $false, $false, $false, $true, 1, $false, "33+61=94", 2) | Out-Null
$d.Content.Find.Execute("13+9=22", $true, $false, $false, $false, $false, $true, 1, $false, "17+19=36", 2) | Out-Null
$d.Content.Find.Execute("88-63=25", $true, $false, $false, $false, $false, $true, 1, $false, "76-35=41", 2) | Out-Null
$d.Content.Find.Execute("59-6=53", $true, $false, $false, $false, $false, $true, 1, $false, "78-9=69", 2) | Out-Null
$d.Content.Find.Execute("63-36=27", $true, $false, $false, $false, $false, $true, 1, $false, "12+73=85", 2) | Out-Null
$d.Content.Find.Execute("21+54=75", $true, $false, $false, $false, $false, $true, 1, $false, "51-17=34", 2) | Out-Null
$d.Content.Find.Execute("26+73=99", $true, $false, $false, $false, $false, $true, 1, $false, "86-66=20", 2) | Out-Null
$d.Content.Find.Execute("52-37=15", $true, $false, $false, $false, $false, $true, 1, $false, "6+49=55", 2) | Out-Null
$d.Content.Find.Execute("37+19=56", $true, $false, $false, $false, $false, $true, 1, $false, "16-1=15", 2) | Out-Null
$d.Content.Find.Execute("84-15=69", $true, $false, $false, $false, $false, $true, 1, $false, "49-44=5", 2) | Out-Null
$d.Content.Find.Execute("83-68=15", $true, $false, $false, $false, $false, $true, 1, $false, "42+35=77", 2) | Out-Null
$d.Content.Find.Execute("5+33=38", $true, $false, $false, $false, $false, $true, 1, $false, "6+75=81", 2) | Out-Null
$d.Content.Find.Execute("97-5=92", $true, $false, $false, $false, $false, $true, 1, $false, "92-24=68", 2) | Out-Null
$d.Content.Find.Execute("11+44=55", $true, $false, $false, $false, $false, $true, 1, $false, "20-8=12", 2) | Out-Null
$d.Content.Find.Execute("30-26=4", $true, $false, $false, $false, $false, $true, 1, $false, "32+37=69", 2) | Out-Null
$d.Content.Find.Execute("66-52=14", $true, $false, $false, $false, $false, $true, 1, $false, "8+17=25", 2) | Out-Null
$d.Content.Find.Execute("1+67=68", $true, $false, $false, $false, $false, $true, 1, $false, "10+20=30", 2) | Out-Null
$d.Content.Find.Execute("22-17=5", $true, $false, $false, $false, $false, $true, 1, $false, "20-10=10", 2) | Out-Null
$d.Content.Find.Execute("37+20=57", $true, $false, $false, $false, $false, $true, 1, $false, "51-42=9", 2) | Out-Null
$d.Content.Find.Execute("4+88=92", $true, $false, $false, $false, $false, $true, 1, $false, "2+61=63", 2) | Out-Null
$d.Content.Find.Execute("26+42=68", $true, $false, $false, $false, $false, $true, 1, $false, "54+12=66", 2) | Out-Null
$d.Content.Find.Execute("57-47=10", $true, $false, $false, $false, $false, $true, 1, $false, "48-10=38", 2) | Out-Null
$d.Content.Find.Execute("14-7=7", $true, $false, $false, $false, $false, $true, 1, $false, "20+49=69", 2) | Out-Null
$d.Content.Find.Execute("11-0=11", $true, $false, $false, $false, $false, $true, 1, $false, "11+63=74", 2) | Out-Null
$d.Content.Find.Execute("82-5=77", $true, $false, $false, $false, $false, $true, 1, $false, "12+41=53", 2) | Out-Null
$d.Content.Find.Execute("91-66=25", $true, $false, $false, $false, $false, $true, 1, $false, "64+19=83", 2) | Out-Null
$d.Content.Find.Execute("66-44=22", $true, $false, $false, $false, $false, $true, 1, $false, "97-27=70", 2) | Out-Null
$d.Content.Find.Execute("36+31=67", $true, $false, $false, $false, $false, $true, 1, $false, "98-24=74", 2) | Out-Null
$d.Content.Find.Execute("83-67=16", $true, $false, $false, $false, $false, $true, 1, $false, "78-31=47", 2) | Out-Null
$d.Content.Find.Execute("66-4=62", $true, $false, $false, $false, $false, $true, 1, $false, "91-8=83", 2) | Out-Null
$d.Content.Find.Execute("65+11=76", $true, $false, $false, $false, $false, $true, 1, $false, "1+92=93", 2) | Out-Null
$d.Content.Find.Execute("37+32=69", $true, $false, $false, $false, $false, $true, 1, $false, "64-34=30", 2) | Out-Null
$d.Content.Find.Execute("51-5=46", $true, $false, $false, $false, $false, $true, 1, $false, "14+79=93", 2) | Out-Null
$d.Content.Find.Execute("91-52=39", $true, $false, $false, $false, $false, $true, 1, $false, "96-82=14", 2) | Out-Null
$d.Content.Find.Execute("33+37=70", $true, $false, $false, $false, $false, $true, 1, $false, "61-36=25", 2) | Out-Null
$d.Content.Find.Execute("57-7=50", $true, $false, $false, $false, $false, $true, 1, $false, "67-23=44", 2) | Out-Null
$d.Content.Find.Execute("1+61=62", $true, $false, $false, $false, $false, $true, 1, $false, "13+61=74", 2) | Out-Null
$d.Content.Find.Execute("32+60=92", $true, $false, $false, $false, $false, $true, 1, $false, "23+45=68", 2) | Out-Null
$d.Content.Find.Execute("39+43=82", $true, $false, $false, $false, $false, $true, 1, $false, "83+2=85", 2) | Out-Null
$d.Content.Find.Execute("22+58=80", $true, $false, $false, $false, $false, $true, 1, $false, "46-4=42", 2) | Out-Null
$d.Content.Find.Execute("94-6=88", $true, $false, $false, $false, $false, $true, 1, $false, "84-43=41", 2) | Out-Null
$d.Content.Find.Execute("37+25=62", $true, $false, $false, $false, $false, $true, 1, $false, "4+66=70", 2) | Out-Null
$d.Content.Find.Execute("64-20=44", $true, $false, $false, $false, $false, $true, 1, $false, "47-11=36", 2) | Out-Null
$d.Content.Find.Execute("32+4=36", $true, $false, $false, $false, $false, $true, 1, $false, "49+8=57", 2) | Out-Null
$d.Content.Find.Execute("13+69=82", $true, $false, $false, $false, $false, $true, 1, $false, "66-49=17", 2) | Out-Null
$d.Content.Find.Execute("45-8=37", $true, $false, $false, $false, $false, $true, 1, $false, "88-22=66", 2) | Out-Null
$d.Content.Find.Execute("72-52=20", $true, $false, $false, $false, $false, $true, 1, $false, "63-43=20", 2) | Out-Null
$d.Content.Find.Execute("29+8=37", $true, $false, $false, $false, $false, $true, 1, $false, "56+5=61", 2) | Out-Null
$d.Content.Find.Execute("38+31=69", $true, $false, $false, $false, $false, $true, 1, $false, "25-7=18", 2) | Out-Null
$d.Content.Find.Execute("75-25=50", $true, $false, $false, $false, $false, $true, 1, $false, "14+33=47", 2) | Out-Null
$d.Content.Find.Execute("46-25=21", $true, $false, $false, $false, $false, $true, 1, $false, "80-13=67", 2) | Out-Null
$d.Content.Find.Execute("87-75=12", $true, $false, $false, $false, $false, $true, 1, $false, "26+38=64", 2) | Out-Null
$d.Content.Find.Execute("50-7=43", $true, $false, $false, $false, $false, $true, 1, $false, "79-58=21", 2) | Out-Null
$d.Content.Find.Execute("20+51=71", $true, $false, $false, $false, $false, $true, 1, $false, "14+58=72", 2) | Out-Null
$d.Content.Find.Execute("28-9=19", $true, $false, $false, $false, $false, $true, 1, $false, "98-78=20", 2) | Out-Null
$d.Content.Find.Execute("78-39=39", $true, $false, $false, $false, $false, $true, 1, $false, "32-15=17", 2) | Out-Null
$d.Content.Find.Execute("0+27=27", $true, $false, $false, $false, $false, $true, 1, $false, "16+9=25", 2) | Out-Null
$d.Content.Find.Execute("97-76=21", $true, $false, $false, $false, $false, $true, 1, $false, "26+70=96", 2) | Out-Null
$d.Content.Find.Execute("80+11=91", $true, $false, $false, $false, $false, $true, 1, $false, "91-0=91", 2) | Out-Null
$d.Content.Find.Execute("61-5=56", $true, $false, $false, $false, $false, $true, 1, $false, "22-16=6", 2) | Out-Null
$d.Content.Find.Execute("80-71=9", $true, $false, $false, $false, $false, $true, 1, $false, "36+40=76", 2) | Out-Null
$d.Content.Find.Execute("14+1=15", $true, $false, $false, $false, $false, $true, 1, $false, "5+6=11", 2) | Out-Null
$d.Content.Find.Execute("68-35=33", $true, $false, $false, $false, $false, $true, 1, $false, "69-57=12", 2) | Out-Null
$d.Content.Find.Execute("60+9=69", $true, $false, $false, $false, $false, $true, 1, $false, "76-69=7", 2) | Out-Null
$d.Content.Find.Execute("78-3=75", $true, $false, $false, $false, $false, $true, 1, $false, "26+58=84", 2) | Out-Null
$d.Content.Find.Execute("23-19=4", $true, $false, $false, $false, $false, $true, 1, $false, "84-55=29", 2) | Out-Null
$d.Content.Find.Execute("22+57=79", $true, $false, $false, $false, $false, $true, 1, $false, "75+10=85", 2) | Out-Null
$d.Content.Find.Execute("48+10=58", $true, $false, $false, $false, $false, $true, 1, $false, "50-35=15", 2) | Out-Null
$d.Content.Find.Execute("46-12=34", $true, $false, $false, $false, $false, $true, 1, $false, "98-73=25", 2) | Out-Null
$d.Content.Find.Execute("83-69=14", $true, $false, $false, $false, $false, $true, 1, $false, "30-20=10", 2) | Out-Null
$d.Content.Find.Execute("29+44=73", $true, $false, $false, $false, $false, $true, 1, $false, "36+56=92", 2) | Out-Null
$d.Content.Find.Execute("55-18=37", $true, $false, $false, $false, $false, $true, 1, $false, "90+6=96", 2) | Out-Null
$d.Content.Find.Execute("60-14=46", $true, $false, $false, $false, $false, $true, 1, $false, "2+17=19", 2) | Out-Null
$d.Content.Find.Execute("96-6=90", $true, $false, $false, $false, $false, $true, 1, $false, "92-35=57", 2) | Out-Null
$d.Content.Find.Execute("39+36=75", $true, $false, $false, $false, $false, $true, 1, $false, "99-48=51", 2) | Out-Null
$d.Content.Find.Execute("98-47=51", $true, $false, $false, $false, $false, $true, 1, $false, "98-81=17", 2) | Out-Null
$d.Content.Find.Execute("56-25=31", $true, $false, $false, $false, $false, $true, 1, $false, "66+7=73", 2) | Out-Null
$d.Content.Find.Execute("96-34=62", $true, $false, $false, $false, $false, $true, 1, $false, "52+16=68", 2) | Out-Null
$d.Content.Find.Execute("53+16=69", $true, $false, $false, $false, $false, $true, 1, $false, "64+29=93", 2) | Out-Null
$d.Content.Find.Execute("61+19=80", $true, $false, $false, $false, $false, $true, 1, $false, "47+30=77", 2) | Out-Null
$d.Content.Find.Execute("70-22=48", $true, $false, $false, $false, $false, $true, 1, $false, "19+76=95", 2) | Out-Null
$d.Content.Find.Execute("75-6=69", $true, $false, $false, $false, $false, $true, 1, $false, "65+4=69", 2) | Out-Null
$d.Content.Find.Execute("6+58=64", $true, $false, $false, $false, $false, $true, 1, $false, "30+22=52", 2) | Out-Null
$d.Content.Find.Execute("66+24=90", $true, $false, $false, $false, $false, $true, 1, $false, "58-11=47", 2) | Out-Null
$d.Content.Find.Execute("2+80=82", $true, $false, $false, $false, $false, $true, 1, $false, "38+9=47", 2) | Out-Null
$d.Content.Find.Execute("48-29=19", $true, $false, $false, $false, $false, $true, 1, $false, "28+25=53", 2) | Out-Null
$d.Content.Find.Execute("1+52=53", $true, $false, $false, $false, $false, $true, 1, $false, "6+54=60", 2) | Out-Null
$d.Content.Find.Execute("5+32=37", $true, $false, $false, $false, $false, $true, 1, $false, "64+35=99", 2) | Out-Null
$d.Content.Find.Execute("43+12=55", $true, $false, $false, $false, $false, $true, 1, $false, "60-39=21", 2) | Out-Null
$d.Content.Find.Execute("26+61=87", $true, $false, $false, $false, $false, $true, 1, $false, "88-31=57", 2) | Out-Null
$d.Content.Find.Execute("51-24=27", $true, $false, $false, $false, $false, $true, 1, $false, "73-22=51", 2) | Out-Null
$d.Content.Find.Execute("18+44=62", $true, $false, $false, $false, $false, $true, 1, $false, "6+10=16", 2) | Out-Null
$d.Content.Find.Execute("92-90=2", $true, $false, $false, $false, $false, $true, 1, $false, "41+20=61", 2) | Out-Null
$d.Content.Find.Execute("18-14=4", $true, $false, $false, $false, $false, $true, 1, $false, "69-22=47", 2) | Out-Null
$d.Content.Find.Execute("14+17=31", $true, $false, $false, $false, $false, $true, 1, $false, "4+38=42", 2) | Out-Null
$d.Content.Find.Execute("50+26=76", $true, $false, $false, $false, $false, $true, 1, $false, "55-46=9", 2) | Out-Null
$d.Content.Find.Execute("75-46=29", $true, $false, $false, $false, $false, $true, 1, $false, "66-3=63", 2) | Out-Null
